$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D2" = "318.63"
    "E2" = "2.98%"
    "D3" = "41.34"
    "E3" = "0.97%"
    "D4" = "5.235"
    "E4" = "2.24%"
    "D5" = "0.07696"
    "E5" = "0.87%"
    "D6" = "1.686"
    "E6" = "5.12%"
    "D7" = "0.9374"
    "E7" = "3.00%"
    "D9" = "0.1279"
    "E9" = "1.23%"
    "D10" = "0.1843"
    "E10" = "1.99%"
    "D11" = "0.09101"
    "E11" = "-0.41%"
    "D12" = "0.04120"
    "E12" = "-5.13%"
    "E13" = "0.23%"
    "D14" = "0.001299"
    "E14" = "4.17%"
    "D15" = "0.006054"
    "E15" = "4.35%"
    "D17" = "3.346"
    "E17" = "-0.12%"
    "D18" = "4.348"
    "E18" = "1.50%"
    "E19" = "1.02%"
    "D20" = "8.414"
    "E20" = "21.21%"
    "E21" = "-2.59%"
    "D22" = "0.2872"
    "E22" = "4.88%"
    "D23" = "0.04054"
    "E23" = "0.19%"
    "D24" = "0.001276"
    "E24" = "0.36%"
    "D25" = "0.004109"
    "E25" = "0.10%"
    "D38" = "0.02534"
    "E38" = "4.54%"
    "D39" = "0.05287"
    "E39" = "0.88%"
    "D40" = "0.007782"
    "E40" = "-0.66%"
    "D41" = "0.1307"
    "E41" = "0.34%"
    "D42" = "0.007067"
    "E42" = "3.78%"
    "D43" = "0.002161"
    "E43" = "16.01%"
    "D44" = "0.008310"
    "E44" = "11.51%"
    "D46" = "0.00006680"
    "E46" = "-2.73%"
    "D47" = "0.00000000754"
    "E47" = "0.41%"
    "D48" = "0.1957"
    "E48" = "35.70%"
    "D49" = "0.004220"
    "E49" = "40.55%"
    "D50" = "0.00002111"
    "E50" = "0.41%"
    "D51" = "0.0002011"
    "E51" = "0.41%"
}

foreach ($cellRef in $changes.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$cellRef]
}
